$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert the two new "emergencyContact" / "emergencyNumber" field
#    paragraphs (plus a fresh blank spacer paragraph) right after the
#    "bloodType" paragraph and before the blank paragraph that currently
#    precedes the "LOG" heading. Word's Range.InsertXML replaces the
#    paragraph that the collapsed range sits at the *start* of, so handing
#    it three <w:p> elements turns that one empty paragraph into three.
# ---------------------------------------------------------------------------

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text
    if ($txt -like "*LOG*") {
        $target = $d.Paragraphs.Item($i - 1)
        break
    }
}

$insertRange = $target.Range
$insertRange.Collapse(1)

$newFieldsXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>emergencyContact</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> varchar(30)</w:t></w:r></w:p>
<w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>emergencyNumber</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> : </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:r><w:t>1</w:t></w:r><w:r><w:t>0)</w:t></w:r></w:p>
<w:p><w:pPr><w:ind w:left="720"/></w:pPr></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertRange.InsertXML($newFieldsXml)

# ---------------------------------------------------------------------------
# 2. Relocate the "_GoBack" bookmark: it used to sit inside the "ID : int(6)"
#    paragraph (between "int" and "("); it now belongs inside the new
#    "emergencyNumber : int(10)" paragraph, between the "1" and the "0)".
# ---------------------------------------------------------------------------

$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

$emergencyNumberPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*emergencyNumber*") {
        $emergencyNumberPara = $p
        break
    }
}

$pStart = $emergencyNumberPara.Range.Start
$pText = $emergencyNumberPara.Range.Text
$offset = $pText.IndexOf("int(1") + 5
$bmPos = $pStart + $offset
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# 3. Drop the stale <w:lastRenderedPageBreak/> marker that sat in front of
#    "timeTaken" - re-running Find & Replace over that exact word rebuilds
#    its run from scratch (keeping the surrounding proofErr markers and
#    paragraph intact) without the page-break marker.
# ---------------------------------------------------------------------------

$rng = $d.Range(0, $d.Content.End)
$rng.Find.Execute("timeTaken", $false, $false, $false, $false, $false, $true, 1, $false, "timeTaken", 2)
